$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '43.676.60'
$ws.Cells.Item(2, 5).Value = '  -0.83%  '
# Row 3
$ws.Cells.Item(3, 4).Value = '2.342.50'
$ws.Cells.Item(3, 5).Value = '  +2.88%  '
# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.09%  '
# Row 5
$ws.Cells.Item(5, 4).Value = "'233.82"
$ws.Cells.Item(5, 5).Value = '  +1.13%  '
# Row 6
$ws.Cells.Item(6, 4).Value = "'0.650"
$ws.Cells.Item(6, 5).Value = '  +2.23%  '
# Row 7
$ws.Cells.Item(7, 4).Value = "'67.37"
$ws.Cells.Item(7, 5).Value = '  +5.43%  '
# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.07%  '
# Row 9
$ws.Cells.Item(9, 5).Value = '  +2.40%  '
# Row 10
$ws.Cells.Item(10, 4).Value = "'0.0973"
$ws.Cells.Item(10, 5).Value = '  -3.94%  '
# Row 11
$ws.Cells.Item(11, 4).Value = "'56.58"
$ws.Cells.Item(11, 5).Value = '  -0.81%  '
# Row 12
$ws.Cells.Item(12, 4).Value = "'27.22"
$ws.Cells.Item(12, 5).Value = '  +1.69%  '
# Row 13
$ws.Cells.Item(13, 4).Value = '2.688.27'
$ws.Cells.Item(13, 5).Value = '  +2.80%  '
# Row 14
$ws.Cells.Item(14, 5).Value = '  -1.27%  '
# Row 15
$ws.Cells.Item(15, 4).Value = "'15.61"
$ws.Cells.Item(15, 5).Value = '  -0.72%  '
# Row 16
$ws.Cells.Item(16, 5).Value = '  +1.64%  '
# Row 17
$ws.Cells.Item(17, 5).Value = '  +1.38%  '
# Row 18
$ws.Cells.Item(18, 4).Value = '2.345.58'
$ws.Cells.Item(18, 5).Value = '  +2.97%  '
# Row 19
$ws.Cells.Item(19, 4).Value = '43.612.38'
$ws.Cells.Item(19, 5).Value = '  -0.80%  '
# Row 20
$ws.Cells.Item(20, 4).Value = '0.0₃0982'
$ws.Cells.Item(20, 5).Value = '  -2.32%  '
# Row 21
$ws.Cells.Item(21, 4).Value = "'74.39"
$ws.Cells.Item(21, 5).Value = '  +0.78%  '
# Row 22
$ws.Cells.Item(22, 4).Value = "'6.29"
$ws.Cells.Item(22, 5).Value = '  +2.72%  '
# Row 23
$ws.Cells.Item(23, 4).Value = "'249.82"
$ws.Cells.Item(23, 5).Value = '  -1.50%  '
# Row 24
$ws.Cells.Item(24, 5).Value = '  +13.46%  '
# Row 25
$ws.Cells.Item(25, 4).Value = "'1.00"
$ws.Cells.Item(25, 5).Value = '  -0.09%  '
# Row 26
$ws.Cells.Item(26, 5).Value = '  -0.75%  '
# Row 27
$ws.Cells.Item(27, 5).Value = '  +0.25%  '
# Row 28
$ws.Cells.Item(28, 4).Value = "'10.01"
$ws.Cells.Item(28, 5).Value = '  -0.42%  '
# Row 29
$ws.Cells.Item(29, 2).Value = 'EthereumClassic'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(29, 4).Value = "'22.25"
$ws.Cells.Item(29, 5).Value = '  +6.58%  '
# Row 30
$ws.Cells.Item(30, 2).Value = 'Monero'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(30, 4).Value = "'174.75"
$ws.Cells.Item(30, 5).Value = '  +1.72%  '
# Row 31
$ws.Cells.Item(31, 4).Value = "'1.48"
$ws.Cells.Item(31, 5).Value = '  +6.21%  '
# Row 32
$ws.Cells.Item(32, 4).Value = "'0.129"
$ws.Cells.Item(32, 5).Value = '  -7.02%  '
# Row 33
$ws.Cells.Item(33, 5).Value = '  +0.80%  '
# Row 34
$ws.Cells.Item(34, 4).Value = "'5.03"
$ws.Cells.Item(34, 5).Value = '  +4.18%  '
# Row 35
$ws.Cells.Item(35, 4).Value = "'0.0691"
$ws.Cells.Item(35, 5).Value = '  -1.40%  '
# Row 36
$ws.Cells.Item(36, 4).Value = "'5.00"
$ws.Cells.Item(36, 5).Value = '  +1.86%  '
# Row 37
$ws.Cells.Item(37, 4).Value = "'2.55"
$ws.Cells.Item(37, 5).Value = '  +9.53%  '
# Row 38
$ws.Cells.Item(38, 4).Value = "'6.57"
$ws.Cells.Item(38, 5).Value = '  +0.28%  '
# Row 39
$ws.Cells.Item(39, 4).Value = "'3.59"
$ws.Cells.Item(39, 5).Value = '  -5.36%  '
# Row 40
$ws.Cells.Item(40, 4).Value = "'0.0255"
$ws.Cells.Item(40, 5).Value = '  -1.83%  '
# Row 41
$ws.Cells.Item(41, 4).Value = "'9.07"
$ws.Cells.Item(41, 5).Value = '  +9.56%  '
# Row 42
$ws.Cells.Item(42, 5).Value = '  +0.00%  '
# Row 43
$ws.Cells.Item(43, 4).Value = "'18.37"
$ws.Cells.Item(43, 5).Value = '  +3.88%  '
# Row 44
$ws.Cells.Item(44, 5).Value = '  +9.05%  '
# Row 45
$ws.Cells.Item(45, 4).Value = "'99.71"
$ws.Cells.Item(45, 5).Value = '  +1.14%  '
# Row 46
$ws.Cells.Item(46, 4).Value = "'1.21"
$ws.Cells.Item(46, 5).Value = '  +0.29%  '
# Row 47
$ws.Cells.Item(47, 4).Value = "'0.0952"
$ws.Cells.Item(47, 5).Value = '  -2.33%  '
# Row 48
$ws.Cells.Item(48, 5).Value = '  +0.33%  '
# Row 49
$ws.Cells.Item(49, 2).Value = 'Celestia'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Cells.Item(49, 4).Value = "'10.09"
$ws.Cells.Item(49, 5).Value = '  -3.31%  '
# Row 50
$ws.Cells.Item(50, 2).Value = 'Maker'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(50, 4).Value = '1.449.14'
$ws.Cells.Item(50, 5).Value = '  -0.14%  '
# Row 51
$ws.Cells.Item(51, 2).Value = 'NEARProtocol'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(51, 4).Value = "'2.32"
$ws.Cells.Item(51, 5).Value = '  +0.82%  '
